$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original Text cell-type: set an explicit
# "Text" number format before writing the values (otherwise Excel auto-converts
# strings like "585.23" or "6.71" into numbers), then clear the format again so
# the cells end up back at the default (unstyled) appearance, matching the source.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '60.761.85'
$ws.Range("E2").Value = '  -3.52%  '
$ws.Range("D3").Value = '2.902.64'
$ws.Range("E3").Value = '  -4.24%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '585.23'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '144.67'
$ws.Range("E6").Value = '  -5.96%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -2.72%  '
$ws.Range("D9").Value = '2.902.57'
$ws.Range("E9").Value = '  -4.16%  '
$ws.Range("E10").Value = '  -2.66%  '
$ws.Range("E11").Value = '  -4.71%  '
$ws.Range("E12").Value = '  -3.71%  '
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("D14").Value = '33.43'
$ws.Range("E14").Value = '  -6.35%  '
$ws.Range("E15").Value = '  +1.64%  '
$ws.Range("D16").Value = '3.382.90'
$ws.Range("E16").Value = '  -4.22%  '
$ws.Range("D17").Value = '60.693.76'
$ws.Range("E17").Value = '  -3.51%  '
$ws.Range("D18").Value = '6.71'
$ws.Range("E18").Value = '  -5.27%  '
$ws.Range("D19").Value = '2.902.60'
$ws.Range("E19").Value = '  -4.30%  '
$ws.Range("D20").Value = '428.38'
$ws.Range("E20").Value = '  -5.33%  '
$ws.Range("D21").Value = '13.56'
$ws.Range("E21").Value = '  -4.81%  '
$ws.Range("D22").Value = '0.681'
$ws.Range("E22").Value = '  -2.42%  '
$ws.Range("E23").Value = '  -5.36%  '
$ws.Range("D24").Value = '80.61'
$ws.Range("E24").Value = '  -2.99%  '
$ws.Range("D25").Value = '10.96'
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").Value = '2.24'
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("D27").Value = '11.91'
$ws.Range("E27").Value = '  -4.25%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").Value = '7.26'
$ws.Range("E30").Value = '  -2.86%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '2.18'
$ws.Range("E31").Value = '  -3.25%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '2.61'
$ws.Range("E32").Value = '  -3.32%  '
$ws.Range("D33").Value = '26.44'
$ws.Range("E33").Value = '  -4.06%  '
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("D35").Value = '0.0₃0876'
$ws.Range("E35").Value = '  +1.87%  '
$ws.Range("E36").Value = '  -3.04%  '
$ws.Range("E37").Value = '  -5.32%  '
$ws.Range("D38").Value = '3.03'
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("E39").Value = '  -3.43%  '
$ws.Range("D40").Value = '49.51'
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("D41").Value = '2.00'
$ws.Range("E41").Value = '  -4.49%  '
$ws.Range("E42").Value = '  -5.75%  '
$ws.Range("D43").Value = '0.298'
$ws.Range("E43").Value = '  -2.93%  '
$ws.Range("D44").Value = '41.33'
$ws.Range("E44").Value = '  -5.54%  '
$ws.Range("D45").Value = '378.60'
$ws.Range("E45").Value = '  -2.93%  '
$ws.Range("D46").Value = '0.0351'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("D47").Value = '2.694.37'
$ws.Range("E47").Value = '  -1.00%  '
$ws.Range("D48").Value = '132.53'
$ws.Range("E48").Value = '  -0.71%  '
$ws.Range("D50").Value = '24.23'
$ws.Range("E50").Value = '  -3.02%  '
$ws.Range("E51").Value = '  -2.51%  '

$priceRange.ClearFormats()

